$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Overview")
$ws.Range("D7").Value = "2016-31-20 04:31:28"

$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("E7").Value = "2016-03-20 04:31:25"

$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("E7").Value = "2016-03-20 04:31:28"
